# Applies the "Formatting and job descriptions" edits to Vincent Chov's resume.
$d = $word.ActiveDocument

function Replace-Text($find, $replace, [bool]$wholeWord = $false) {
    $d.Content.Find.Execute($find, $true, $wholeWord, $false, $false, $false,
                             $true, 1, $false, $replace, 2) | Out-Null
}

# Phone number update
Replace-Text "(203) 543-3274 " "(203) 307-5485 "

# Education section: title-case school names
Replace-Text "UNIVERSITY OF CONNECTICUT, Storrs, CT" "University of Connecticut, Storrs, CT"
Replace-Text "HARVARD UNIVERSITY, EXTENSION SCHOOL, Cambridge, MA" "Harvard University, Extension School, Cambridge, MA"

# Technical skills
Replace-Text "SQL, Docker, Flask, React, Linux, Git, Sass, Selenium, Bootstrap" "Microsoft SQL Server, Docker, Flask, React, Linux, Git, Sass, Selenium, Bootstrap"

# Professional experience: NGP VAN role
Replace-Text "SOFTWARE ENGINEERING INTERN (FULL-TIME)" "Software Engineering Intern (Full-Time)"
Replace-Text "Enhanced legacy codebase to enable political campaigns to organize volunteers and fundraise." "Enhanced legacy ASP.NET codebase to enable political campaigns to organize volunteers and fundraise."
Replace-Text "Wrote Selenium unit tests for use in a Continuous Integration / Continuous Delivery (CI/CD) server." "Ensured correctness and quality of ETL processes using Selenium unit tests."

# Professional experience: second role
Replace-Text "SOFTWARE ENGINEER" "Software Engineer" $true
Replace-Text "Performed database migrations to use SQL Server and created a temporary Python/Flask front-end." "Developed an ETL-like process to enable a switch to using Microsoft SQL Server."
Replace-Text "Deployed a web application for internal use on a Gunicorn server and Nginx reverse proxy, hosted on a Debian VPS." "Deployed a Flask API for internal use on a Gunicorn server and Nginx reverse proxy."

# Section heading: PROJECT/VOLUNTEER EXPERIENCE -> OPEN SOURCE/VOLUNTEER EXPERIENCE
Replace-Text "PROJECT/VOLUNTEER EXPERIENCE" "OPEN SOURCE/VOLUNTEER EXPERIENCE"

# Open source / volunteer experience entry
Replace-Text "FULL STACK DEVELOPER" "Full Stack Developer"
Replace-Text "JANUARY 2017 - PRESENT" "JANUARY 2017 - AUGUST 2019"
Replace-Text "Improved Swing dance event coordination process by writing a web crawler with Selenium." "Automated portions of the event coordination process by writing a Selenium-based web crawler."
Replace-Text "Introduced PHP templating to maintain consistency of information across the website." "Altered deployment scheme to enable the use of Git in management of configuration settings."
Replace-Text "Made their website mobile-friendly using responsive design." "Leveraged the Understrap framework to allow for the use of Bootstrap in WordPress."
